$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 6464.1875  # H28: 6435.8125 -> 6464.1875
$ws.Cells.Item(28, 9).Value = 6863.2  # I28: 6435.8125 -> 6863.2
$ws.Cells.Item(28, 10).Value = 479  # J28: 0 -> 479
$ws.Cells.Item(28, 11).Value = 6863.2  # K28: 6435.8125 -> 6863.2
$ws.Cells.Item(28, 12).Value = 479  # L28: 0 -> 479
$ws.Cells.Item(28, 13).Value = -6378.2  # M28: -5950.8125 -> -6378.2
$ws.Cells.Item(28, 14).Value = -1449  # N28: (empty) -> -1449
$ws.Cells.Item(87, 8).Value = 0  # H87: 75000 -> 0
$ws.Cells.Item(87, 10).Value = 0  # J87: 75000 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 75000 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: -77496 -> (empty)
$ws.Cells.Item(90, 8).Value = 0  # H90: 75000 -> 0
$ws.Cells.Item(90, 10).Value = 0  # J90: 75000 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 225000 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: -237480 -> (empty)
$ws.Cells.Item(100, 8).Value = 5137.25  # H100: 5305.5884 -> 5137.25
$ws.Cells.Item(100, 10).Value = 8249.875  # J100: 8222 -> 8249.875
$ws.Cells.Item(100, 12).Value = 8249.875  # L100: 8222 -> 8249.875
$ws.Cells.Item(100, 14).Value = -9331.875  # N100: -9304 -> -9331.875
$ws.Cells.Item(107, 8).Value = 914.1111  # H107: 460.95 -> 914.1111
$ws.Cells.Item(107, 9).Value = 914.1111  # I107: 502.83334 -> 914.1111
$ws.Cells.Item(107, 10).Value = 0  # J107: 84 -> 0
$ws.Cells.Item(107, 11).Value = 914.1111  # K107: 502.83334 -> 914.1111
$ws.Cells.Item(107, 12).Value = 0  # L107: 84 -> 0
$ws.Cells.Item(107, 13).Value = 1005.8889  # M107: 1417.16666 -> 1005.8889
$ws.Cells.Item(107, 14).ClearContents()  # N107: -3924 -> (empty)
$ws.Cells.Item(111, 8).Value = 47620870  # H111: 55557436 -> 47620870
$ws.Cells.Item(111, 9).Value = 1751.8  # I111: 2256.2 -> 1751.8
$ws.Cells.Item(111, 10).Value = 166668670  # J111: 333333340 -> 166668670
$ws.Cells.Item(111, 11).Value = 5255.4  # K111: 6768.599999999999 -> 5255.4
$ws.Cells.Item(111, 12).Value = 500006010  # L111: 1000000020 -> 500006010
$ws.Cells.Item(111, 13).Value = -2188.4  # M111: -3701.599999999999 -> -2188.4
$ws.Cells.Item(111, 14).Value = -500012144  # N111: -1000006154 -> -500012144
$ws.Cells.Item(118, 8).Value = 899.26666  # H118: 1468.0625 -> 899.26666
$ws.Cells.Item(118, 10).Value = 450  # J118: 1814.2858 -> 450
$ws.Cells.Item(118, 12).Value = 1350  # L118: 5442.857400000001 -> 1350
$ws.Cells.Item(118, 14).Value = -4664  # N118: -8756.857400000001 -> -4664
$ws.Cells.Item(136, 8).Value = 79243.75  # H136: 79662.664 -> 79243.75
$ws.Cells.Item(136, 10).Value = 79243.75  # J136: 79662.664 -> 79243.75
$ws.Cells.Item(136, 12).Value = 79243.75  # L136: 79662.664 -> 79243.75
$ws.Cells.Item(136, 14).Value = -89443.75  # N136: -89862.664 -> -89443.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1899.0934  # H32: 1921.1948 -> 1899.0934
$ws.Cells.Item(32, 10).Value = 0  # J32: 2750 -> 0
$ws.Cells.Item(32, 12).Value = 0  # L32: 2750 -> 0
$ws.Cells.Item(32, 14).ClearContents()  # N32: -3324 -> (empty)
$ws.Cells.Item(35, 8).Value = 0  # H35: 2000 -> 0
$ws.Cells.Item(35, 9).Value = 0  # I35: 2000 -> 0
$ws.Cells.Item(35, 11).Value = 0  # K35: 2000 -> 0
$ws.Cells.Item(35, 13).ClearContents()  # M35: -1594 -> (empty)
$ws.Cells.Item(74, 8).Value = 1888.6316  # H74: 1683 -> 1888.6316
$ws.Cells.Item(74, 9).Value = 1882.4445  # I74: 1667.9048 -> 1882.4445
$ws.Cells.Item(74, 11).Value = 1882.4445  # K74: 1667.9048 -> 1882.4445
$ws.Cells.Item(74, 13).Value = -1008.4445  # M74: -793.9048 -> -1008.4445
$ws.Cells.Item(77, 8).Value = 1888.6316  # H77: 1683 -> 1888.6316
$ws.Cells.Item(77, 9).Value = 1882.4445  # I77: 1667.9048 -> 1882.4445
$ws.Cells.Item(77, 11).Value = 9412.2225  # K77: 8339.523999999999 -> 9412.2225
$ws.Cells.Item(77, 13).Value = -5044.2225  # M77: -3971.523999999999 -> -5044.2225
$ws.Cells.Item(132, 8).Value = 1757.7778  # H132: 1759.537 -> 1757.7778
$ws.Cells.Item(132, 9).Value = 1783.8125  # I132: 1785.7916 -> 1783.8125
$ws.Cells.Item(132, 11).Value = 5351.4375  # K132: 5357.3748 -> 5351.4375
$ws.Cells.Item(132, 13).Value = -2821.4375  # M132: -2827.3748 -> -2821.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2454.25  # H134: 2475.739 -> 2454.25
$ws.Cells.Item(134, 9).Value = 2355.1  # I134: 2375.8948 -> 2355.1
$ws.Cells.Item(134, 11).Value = 7065.299999999999  # K134: 7127.6844 -> 7065.299999999999
$ws.Cells.Item(134, 13).Value = -4530.299999999999  # M134: -4592.6844 -> -4530.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2645.7827  # H31: 3021.95 -> 2645.7827
$ws.Cells.Item(31, 9).Value = 1305.091  # I31: 1660.2222 -> 1305.091
$ws.Cells.Item(31, 10).Value = 3874.75  # J31: 4136.091 -> 3874.75
$ws.Cells.Item(31, 11).Value = 1305.091  # K31: 1660.2222 -> 1305.091
$ws.Cells.Item(31, 12).Value = 3874.75  # L31: 4136.091 -> 3874.75
$ws.Cells.Item(31, 13).Value = -1010.091  # M31: -1365.2222 -> -1010.091
$ws.Cells.Item(31, 14).Value = -4464.75  # N31: -4726.091 -> -4464.75
$ws.Cells.Item(34, 8).Value = 2645.7827  # H34: 3021.95 -> 2645.7827
$ws.Cells.Item(34, 9).Value = 1305.091  # I34: 1660.2222 -> 1305.091
$ws.Cells.Item(34, 10).Value = 3874.75  # J34: 4136.091 -> 3874.75
$ws.Cells.Item(34, 11).Value = 1305.091  # K34: 1660.2222 -> 1305.091
$ws.Cells.Item(34, 12).Value = 3874.75  # L34: 4136.091 -> 3874.75
$ws.Cells.Item(34, 13).Value = -1103.091  # M34: -1458.2222 -> -1103.091
$ws.Cells.Item(34, 14).Value = -4278.75  # N34: -4540.091 -> -4278.75
$ws.Cells.Item(58, 8).Value = 3635.2727  # H58: 3676.7 -> 3635.2727
$ws.Cells.Item(58, 9).Value = 2996.75  # I58: 2969.125 -> 2996.75
$ws.Cells.Item(58, 10).Value = 5338  # J58: 6507 -> 5338
$ws.Cells.Item(58, 11).Value = 2996.75  # K58: 2969.125 -> 2996.75
$ws.Cells.Item(58, 12).Value = 5338  # L58: 6507 -> 5338
$ws.Cells.Item(58, 13).Value = -2793.75  # M58: -2766.125 -> -2793.75
$ws.Cells.Item(58, 14).Value = -5744  # N58: -6913 -> -5744
$ws.Cells.Item(94, 8).Value = 4125.7144  # H94: 2420.3076 -> 4125.7144
$ws.Cells.Item(94, 10).Value = 5267  # J94: 2365.2 -> 5267
$ws.Cells.Item(94, 12).Value = 5267  # L94: 2365.2 -> 5267
$ws.Cells.Item(94, 14).Value = -6169  # N94: -3267.2 -> -6169
$ws.Cells.Item(99, 8).Value = 1003741.6  # H99: 1004152.5 -> 1003741.6
$ws.Cells.Item(99, 9).Value = 1115058.1  # I99: 1115514.6 -> 1115058.1
$ws.Cells.Item(99, 11).Value = 1115058.1  # K99: 1115514.6 -> 1115058.1
$ws.Cells.Item(99, 13).Value = -1113560.1  # M99: -1114016.6 -> -1113560.1
$ws.Cells.Item(122, 8).Value = 900.125  # H122: 925.125 -> 900.125
$ws.Cells.Item(122, 9).Value = 720.8889  # I122: 773.5 -> 720.8889
$ws.Cells.Item(122, 10).Value = 1130.5714  # J122: 1076.75 -> 1130.5714
$ws.Cells.Item(122, 11).Value = 2162.6667  # K122: 2320.5 -> 2162.6667
$ws.Cells.Item(122, 12).Value = 3391.7142  # L122: 3230.25 -> 3391.7142
$ws.Cells.Item(122, 13).Value = 287.3332999999998  # M122: 129.5 -> 287.3332999999998
$ws.Cells.Item(122, 14).Value = -8291.7142  # N122: -8130.25 -> -8291.7142
$ws.Cells.Item(126, 8).Value = 1003741.6  # H126: 1004152.5 -> 1003741.6
$ws.Cells.Item(126, 9).Value = 1115058.1  # I126: 1115514.6 -> 1115058.1
$ws.Cells.Item(126, 11).Value = 3345174.3  # K126: 3346543.8 -> 3345174.3
$ws.Cells.Item(126, 13).Value = -3342704.3  # M126: -3344073.8 -> -3342704.3
$ws.Cells.Item(132, 8).Value = 1603.0952  # H132: 1668.3 -> 1603.0952
$ws.Cells.Item(132, 9).Value = 1370.5555  # I132: 1433.5883 -> 1370.5555
$ws.Cells.Item(132, 11).Value = 4111.666499999999  # K132: 4300.7649 -> 4111.666499999999
$ws.Cells.Item(132, 13).Value = -1581.666499999999  # M132: -1770.7649 -> -1581.666499999999
$ws.Cells.Item(134, 8).Value = 5083.9375  # H134: 5077.6875 -> 5083.9375
$ws.Cells.Item(134, 9).Value = 3884.4  # I134: 4149.4443 -> 3884.4
$ws.Cells.Item(134, 10).Value = 7083.1665  # J134: 6271.143 -> 7083.1665
$ws.Cells.Item(134, 11).Value = 11653.2  # K134: 12448.3329 -> 11653.2
$ws.Cells.Item(134, 12).Value = 21249.4995  # L134: 18813.429 -> 21249.4995
$ws.Cells.Item(134, 13).Value = -9118.200000000001  # M134: -9913.332900000001 -> -9118.200000000001
$ws.Cells.Item(134, 14).Value = -26319.4995  # N134: -23883.429 -> -26319.4995
$ws.Cells.Item(136, 8).Value = 3635.2727  # H136: 3676.7 -> 3635.2727
$ws.Cells.Item(136, 9).Value = 2996.75  # I136: 2969.125 -> 2996.75
$ws.Cells.Item(136, 10).Value = 5338  # J136: 6507 -> 5338
$ws.Cells.Item(136, 11).Value = 8990.25  # K136: 8907.375 -> 8990.25
$ws.Cells.Item(136, 12).Value = 16014  # L136: 19521 -> 16014
$ws.Cells.Item(136, 13).Value = -6440.25  # M136: -6357.375 -> -6440.25
$ws.Cells.Item(136, 14).Value = -21114  # N136: -24621 -> -21114
$ws.Cells.Item(141, 8).Value = 149991.67  # H141: 216551.67 -> 149991.67
$ws.Cells.Item(141, 10).Value = 149991.67  # J141: 216551.67 -> 149991.67
$ws.Cells.Item(141, 12).Value = 149991.67  # L141: 216551.67 -> 149991.67
$ws.Cells.Item(141, 14).Value = -160351.67  # N141: -226911.67 -> -160351.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 3370.7273  # H107: 2907.4 -> 3370.7273
$ws.Cells.Item(107, 9).Value = 6400.2  # I107: 4899.2 -> 6400.2
$ws.Cells.Item(107, 10).Value = 846.1667  # J107: 915.6 -> 846.1667
$ws.Cells.Item(107, 11).Value = 19200.6  # K107: 14697.6 -> 19200.6
$ws.Cells.Item(107, 12).Value = 2538.5001  # L107: 2746.8 -> 2538.5001
$ws.Cells.Item(107, 13).Value = -17280.6  # M107: -12777.6 -> -17280.6
$ws.Cells.Item(107, 14).Value = -6378.5001  # N107: -6586.8 -> -6378.5001
$ws.Cells.Item(122, 8).Value = 3045.1667  # H122: 3187.6667 -> 3045.1667
$ws.Cells.Item(122, 10).Value = 3222  # J122: 3377.4546 -> 3222
$ws.Cells.Item(122, 12).Value = 28998  # L122: 30397.0914 -> 28998
$ws.Cells.Item(122, 14).Value = -33898  # N122: -35297.0914 -> -33898
$ws.Cells.Item(140, 8).Value = 1734.0834  # H140: 1809.909 -> 1734.0834
$ws.Cells.Item(140, 9).Value = 1734.0834  # I140: 1809.909 -> 1734.0834
$ws.Cells.Item(140, 11).Value = 5202.2502  # K140: 5429.727000000001 -> 5202.2502
$ws.Cells.Item(140, 13).Value = -22.2502000000004  # M140: -249.7270000000008 -> -22.2502000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 17433.666  # H97: 38127.25 -> 17433.666
$ws.Cells.Item(97, 9).Value = 852.25  # I97: 1010 -> 852.25
$ws.Cells.Item(97, 10).Value = 30698.8  # J97: 50499.668 -> 30698.8
$ws.Cells.Item(97, 11).Value = 852.25  # K97: 1010 -> 852.25
$ws.Cells.Item(97, 12).Value = 30698.8  # L97: 50499.668 -> 30698.8
$ws.Cells.Item(97, 13).Value = -356.25  # M97: -514 -> -356.25
$ws.Cells.Item(97, 14).Value = -31690.8  # N97: -51491.668 -> -31690.8
$ws.Cells.Item(134, 8).Value = 46382.715  # H134: 49099.832 -> 46382.715
$ws.Cells.Item(134, 10).Value = 46382.715  # J134: 49099.832 -> 46382.715
$ws.Cells.Item(134, 12).Value = 139148.145  # L134: 147299.496 -> 139148.145
$ws.Cells.Item(134, 14).Value = -144218.145  # N134: -152369.496 -> -144218.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6985.154  # H40: 6957.4614 -> 6985.154
$ws.Cells.Item(40, 9).Value = 5700.5  # I40: 4510.4443 -> 5700.5
$ws.Cells.Item(40, 10).Value = 7556.1113  # J40: 8252.941000000001 -> 7556.1113
$ws.Cells.Item(40, 11).Value = 5700.5  # K40: 4510.4443 -> 5700.5
$ws.Cells.Item(40, 12).Value = 7556.1113  # L40: 8252.941000000001 -> 7556.1113
$ws.Cells.Item(40, 13).Value = -5564.5  # M40: -4374.4443 -> -5564.5
$ws.Cells.Item(40, 14).Value = -7828.1113  # N40: -8524.941000000001 -> -7828.1113
$ws.Cells.Item(46, 8).Value = 2054.7  # H46: 2396.8 -> 2054.7
$ws.Cells.Item(46, 9).Value = 1508.6666  # I46: 1778.5714 -> 1508.6666
$ws.Cells.Item(46, 10).Value = 6969  # J46: 3839.3333 -> 6969
$ws.Cells.Item(46, 11).Value = 1508.6666  # K46: 1778.5714 -> 1508.6666
$ws.Cells.Item(46, 12).Value = 6969  # L46: 3839.3333 -> 6969
$ws.Cells.Item(46, 13).Value = -1320.6666  # M46: -1590.5714 -> -1320.6666
$ws.Cells.Item(46, 14).Value = -7345  # N46: -4215.3333 -> -7345
$ws.Cells.Item(135, 8).Value = 40000  # H135: 0 -> 40000
$ws.Cells.Item(135, 10).Value = 40000  # J135: 0 -> 40000
$ws.Cells.Item(135, 12).Value = 40000  # L135: 0 -> 40000
$ws.Cells.Item(135, 14).Value = -50140  # N135: (empty) -> -50140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 11679.5  # H41: 11259.444 -> 11679.5
$ws.Cells.Item(41, 10).Value = 11299.667  # J41: 10449.5 -> 11299.667
$ws.Cells.Item(41, 12).Value = 11299.667  # L41: 10449.5 -> 11299.667
$ws.Cells.Item(41, 14).Value = -12079.667  # N41: -11229.5 -> -12079.667
$ws.Cells.Item(46, 8).Value = 0  # H46: 83429 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 83429 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 83429 -> 0
$ws.Cells.Item(46, 14).ClearContents()  # N46: -83891 -> (empty)
$ws.Cells.Item(132, 8).Value = 2229.2  # H132: 2332.4443 -> 2229.2
$ws.Cells.Item(132, 9).Value = 2255.2222  # I132: 2374.625 -> 2255.2222
$ws.Cells.Item(132, 11).Value = 6765.6666  # K132: 7123.875 -> 6765.6666
$ws.Cells.Item(132, 13).Value = -4235.6666  # M132: -4593.875 -> -4235.6666
$ws.Cells.Item(134, 8).Value = 0  # H134: 83429 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134: 83429 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134: 250287 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: -255357 -> (empty)
$ws.Cells.Item(136, 8).Value = 1889.8334  # H136: 1656.5 -> 1889.8334
$ws.Cells.Item(136, 9).Value = 1084.75  # I136: 1042 -> 1084.75
$ws.Cells.Item(136, 11).Value = 3254.25  # K136: 3126 -> 3254.25
$ws.Cells.Item(136, 13).Value = -704.25  # M136: -576 -> -704.25
